$wb = $excel.ActiveWorkbook

# ALC row 4
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 68.8
$ws.Range("I4").Value = 68.8
$ws.Range("K4").Value = 68.8
$ws.Range("M4").Value = 45.2

# ALC row 6
$ws.Range("H6").Value = 550
$ws.Range("I6").Value = 216.66667
$ws.Range("J6").Value = 1050
$ws.Range("K6").Value = 650.00001
$ws.Range("L6").Value = 3150
$ws.Range("M6").Value = -538.00001
$ws.Range("N6").Value = -3374

# ALC row 18
$ws.Range("H18").Value = 295
$ws.Range("I18").Value = 295
$ws.Range("K18").Value = 295
$ws.Range("M18").Value = -11

# ALC row 96
$ws.Range("H96").Value = 8500
$ws.Range("I96").Value = 8500
$ws.Range("K96").Value = 25500
$ws.Range("M96").Value = -24127

# ALC row 97
$ws.Range("H97").Value = 1329.3334
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 1329.3334
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 3988.0002
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -4980.0002

# ALC row 99
$ws.Range("H99").Value = 7595.5
$ws.Range("I99").Value = 7595.5
$ws.Range("K99").Value = 22786.5
$ws.Range("M99").Value = -21288.5

# ALC row 100
$ws.Range("H100").Value = 3731.818
$ws.Range("I100").Value = 2200
$ws.Range("K100").Value = 2200
$ws.Range("M100").Value = -1659

# ALC row 101
$ws.Range("H101").Value = 1359
$ws.Range("I101").Value = 374
$ws.Range("J101").Value = 3329
$ws.Range("K101").Value = 1122
$ws.Range("L101").Value = 9987
$ws.Range("M101").Value = 500
$ws.Range("N101").Value = -13231

# ARM row 6
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 150010000
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()

# ARM row 32
$ws.Range("H32").Value = 3029.7122
$ws.Range("I32").Value = 2629.9539
$ws.Range("K32").Value = 2629.9539
$ws.Range("M32").Value = -2342.9539

# ARM row 74
$ws.Range("H74").Value = 30306442
$ws.Range("I74").Value = 55560160
$ws.Range("J74").Value = 1980
$ws.Range("K74").Value = 55560160
$ws.Range("L74").Value = 1980
$ws.Range("M74").Value = -55559286
$ws.Range("N74").Value = -3728

# ARM row 77
$ws.Range("H77").Value = 30306442
$ws.Range("I77").Value = 55560160
$ws.Range("J77").Value = 1980
$ws.Range("K77").Value = 277800800
$ws.Range("L77").Value = 9900
$ws.Range("M77").Value = -277796432
$ws.Range("N77").Value = -18636

# BSM row 53
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 149998
$ws.Range("J53").Value = 149998
$ws.Range("L53").Value = 149998
$ws.Range("N53").Value = -151146

# BSM row 118
$ws.Range("H118").Value = 62000
$ws.Range("J118").Value = 62000
$ws.Range("L118").Value = 62000
$ws.Range("N118").Value = -65314

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5488.9375
$ws.Range("I58").Value = 2711.111
$ws.Range("J58").Value = 9060.429
$ws.Range("K58").Value = 2711.111
$ws.Range("L58").Value = 9060.429
$ws.Range("M58").Value = -2508.111
$ws.Range("N58").Value = -9466.429

# CRP row 62
$ws.Range("H62").Value = 7598.5
$ws.Range("J62").Value = 11750.833
$ws.Range("L62").Value = 11750.833
$ws.Range("N62").Value = -12998.833

# CRP row 65
$ws.Range("H65").Value = 7598.5
$ws.Range("J65").Value = 11750.833
$ws.Range("L65").Value = 58754.165
$ws.Range("N65").Value = -64994.165

# CRP row 94
$ws.Range("H94").Value = 2125.5833
$ws.Range("I94").Value = 1705.75
$ws.Range("J94").Value = 2335.5
$ws.Range("K94").Value = 1705.75
$ws.Range("L94").Value = 2335.5
$ws.Range("M94").Value = -1254.75
$ws.Range("N94").Value = -3237.5

# CRP row 134
$ws.Range("H134").Value = 5903.4287
$ws.Range("I134").Value = 3975
$ws.Range("K134").Value = 11925
$ws.Range("M134").Value = -9390

# CRP row 136
$ws.Range("H136").Value = 5488.9375
$ws.Range("I136").Value = 2711.111
$ws.Range("J136").Value = 9060.429
$ws.Range("K136").Value = 8133.333
$ws.Range("L136").Value = 27181.287
$ws.Range("M136").Value = -5583.333
$ws.Range("N136").Value = -32281.287

# CRP row 140
$ws.Range("H140").Value = 63999.5
$ws.Range("J140").Value = 63999.5
$ws.Range("L140").Value = 63999.5
$ws.Range("N140").Value = -74359.5

# CUL row 16
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 1250
$ws.Range("I16").Value = 750
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 2250
$ws.Range("L16").Value = 6000
$ws.Range("M16").Value = -2077
$ws.Range("N16").Value = -6346

# CUL row 19
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()

# GSM row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 92.38461
$ws.Range("I2").Value = 19.5
$ws.Range("K2").Value = 19.5
$ws.Range("M2").Value = 93.5

# GSM row 5
$ws.Range("H5").Value = 199.2
$ws.Range("I5").Value = 199.2
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 199.2
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -87.19999999999999
$ws.Range("N5").ClearContents()

# GSM row 53
$ws.Range("H53").Value = 43999
$ws.Range("J53").Value = 43999
$ws.Range("L53").Value = 43999
$ws.Range("N53").Value = -45261

# GSM row 99
$ws.Range("H99").Value = 25167.084
$ws.Range("I99").Value = 14073
$ws.Range("J99").Value = 40698.8
$ws.Range("K99").Value = 14073
$ws.Range("L99").Value = 40698.8
$ws.Range("M99").Value = -11827
$ws.Range("N99").Value = -45190.8

# GSM row 132
$ws.Range("H132").Value = 68797
$ws.Range("I132").Value = 128195.75
$ws.Range("J132").Value = 9398.25
$ws.Range("K132").Value = 384587.25
$ws.Range("L132").Value = 28194.75
$ws.Range("M132").Value = -382057.25
$ws.Range("N132").Value = -33254.75

# GSM row 139
$ws.Range("H139").Value = 74884
$ws.Range("J139").Value = 74884
$ws.Range("L139").Value = 74884
$ws.Range("N139").Value = -85164

# WVR row 4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 6999.6665
$ws.Range("J4").Value = 6999.6665
$ws.Range("L4").Value = 6999.6665
$ws.Range("N4").Value = -7225.6665

# WVR row 107
$ws.Range("H107").Value = 468.76923
$ws.Range("I107").Value = 413.0909
$ws.Range("K107").Value = 1239.2727
$ws.Range("M107").Value = 680.7273

# WVR row 136
$ws.Range("H136").Value = 1864.826
$ws.Range("I136").Value = 1864.826
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5594.478
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -3044.478
$ws.Range("N136").ClearContents()
